# Update cached market-price / profit figures across the Leve-crafting
# tables on each job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below were refreshed from the scheduled market-data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3553.8462
$ws.Range("I64").Value = 2971.4285
$ws.Range("J64").Value = 4233.3335
$ws.Range("K64").Value = 2971.4285
$ws.Range("L64").Value = 4233.3335
$ws.Range("M64").Value = -2723.4285
$ws.Range("N64").Value = -4729.3335
$ws.Range("H67").Value = 3553.8462
$ws.Range("I67").Value = 2971.4285
$ws.Range("J67").Value = 4233.3335
$ws.Range("K67").Value = 2971.4285
$ws.Range("L67").Value = 4233.3335
$ws.Range("M67").Value = -2113.4285
$ws.Range("N67").Value = -5949.3335
$ws.Range("H69").Value = 1573.2142
$ws.Range("I69").Value = 2300
$ws.Range("J69").Value = 1517.3077
$ws.Range("K69").Value = 6900
$ws.Range("L69").Value = 4551.9231
$ws.Range("M69").Value = -6026
$ws.Range("N69").Value = -6299.9231
$ws.Range("H72").Value = 1573.2142
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 1517.3077
$ws.Range("K72").Value = 20700
$ws.Range("L72").Value = 13655.7693
$ws.Range("M72").Value = -16332
$ws.Range("N72").Value = -22391.7693
$ws.Range("H116").Value = 5166.222
$ws.Range("I116").Value = 1800
$ws.Range("J116").Value = 5587
$ws.Range("K116").Value = 1800
$ws.Range("L116").Value = 5587
$ws.Range("M116").Value = 1642
$ws.Range("N116").Value = -12471
$ws.Range("H132").Value = 3787.652
$ws.Range("I132").Value = 3914.3635
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 11743.0905
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -9213.0905
$ws.Range("N132").Value = -8060
$ws.Range("H138").Value = 2292.6064
$ws.Range("I138").Value = 1335.3462
$ws.Range("J138").Value = 2658.6177
$ws.Range("K138").Value = 4006.0386
$ws.Range("L138").Value = 7975.853099999999
$ws.Range("M138").Value = 1133.9614
$ws.Range("N138").Value = -18255.8531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15783.488
$ws.Range("I32").Value = 16934.455
$ws.Range("J32").Value = 5936.3335
$ws.Range("K32").Value = 16934.455
$ws.Range("L32").Value = 5936.3335
$ws.Range("M32").Value = -16647.455
$ws.Range("N32").Value = -6510.3335
$ws.Range("H61").Value = 1906.1428
$ws.Range("I61").Value = 1868.2778
$ws.Range("J61").Value = 2133.3333
$ws.Range("K61").Value = 1868.2778
$ws.Range("L61").Value = 2133.3333
$ws.Range("M61").Value = -1656.2778
$ws.Range("N61").Value = -2557.3333
$ws.Range("H74").Value = 27028772
$ws.Range("I74").Value = 27779550
$ws.Range("J74").Value = 800
$ws.Range("K74").Value = 27779550
$ws.Range("L74").Value = 800
$ws.Range("M74").Value = -27778676
$ws.Range("N74").Value = -2548
$ws.Range("H77").Value = 27028772
$ws.Range("I77").Value = 27779550
$ws.Range("J77").Value = 800
$ws.Range("K77").Value = 138897750
$ws.Range("L77").Value = 4000
$ws.Range("M77").Value = -138893382
$ws.Range("N77").Value = -12736
$ws.Range("H97").Value = 1266.5588
$ws.Range("I97").Value = 1185.5416
$ws.Range("J97").Value = 1461
$ws.Range("K97").Value = 1185.5416
$ws.Range("L97").Value = 1461
$ws.Range("M97").Value = -689.5416
$ws.Range("N97").Value = -2453
$ws.Range("H132").Value = 52366.3
$ws.Range("I132").Value = 2546.1052
$ws.Range("J132").Value = 138419.36
$ws.Range("K132").Value = 7638.3156
$ws.Range("L132").Value = 415258.08
$ws.Range("M132").Value = -5108.3156
$ws.Range("N132").Value = -420318.08
$ws.Range("H136").Value = 1906.1428
$ws.Range("I136").Value = 1868.2778
$ws.Range("J136").Value = 2133.3333
$ws.Range("K136").Value = 5604.8334
$ws.Range("L136").Value = 6399.999899999999
$ws.Range("M136").Value = -3054.8334
$ws.Range("N136").Value = -11499.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2586.1667
$ws.Range("I20").Value = 4002.6667
$ws.Range("J20").Value = 1169.6666
$ws.Range("K20").Value = 4002.6667
$ws.Range("L20").Value = 1169.6666
$ws.Range("M20").Value = -3755.6667
$ws.Range("N20").Value = -1663.6666
$ws.Range("H86").Value = 1496.4445
$ws.Range("I86").Value = 1363.8948
$ws.Range("J86").Value = 1811.25
$ws.Range("K86").Value = 1363.8948
$ws.Range("L86").Value = 1811.25
$ws.Range("M86").Value = -240.8948
$ws.Range("N86").Value = -4057.25
$ws.Range("H89").Value = 1496.4445
$ws.Range("I89").Value = 1363.8948
$ws.Range("J89").Value = 1811.25
$ws.Range("K89").Value = 6819.474
$ws.Range("L89").Value = 9056.25
$ws.Range("M89").Value = -1203.474
$ws.Range("N89").Value = -20288.25
$ws.Range("H99").Value = 1728.1818
$ws.Range("I99").Value = 1842
$ws.Range("J99").Value = 1633.3334
$ws.Range("K99").Value = 1842
$ws.Range("L99").Value = 1633.3334
$ws.Range("M99").Value = -344
$ws.Range("N99").Value = -4629.3334
$ws.Range("H134").Value = 23074.02
$ws.Range("I134").Value = 28179.195
$ws.Range("J134").Value = 2142.8
$ws.Range("K134").Value = 84537.58499999999
$ws.Range("L134").Value = 6428.400000000001
$ws.Range("M134").Value = -82002.58499999999
$ws.Range("N134").Value = -11498.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 18522086
$ws.Range("I99").Value = 3228.524
$ws.Range("J99").Value = 83338080
$ws.Range("K99").Value = 3228.524
$ws.Range("L99").Value = 83338080
$ws.Range("M99").Value = -1730.524
$ws.Range("N99").Value = -83341076
$ws.Range("H126").Value = 18522086
$ws.Range("I126").Value = 3228.524
$ws.Range("J126").Value = 83338080
$ws.Range("K126").Value = 9685.572
$ws.Range("L126").Value = 250014240
$ws.Range("M126").Value = -7215.572
$ws.Range("N126").Value = -250019180
$ws.Range("H132").Value = 16394.834
$ws.Range("I132").Value = 20334.148
$ws.Range("J132").Value = 4576.8887
$ws.Range("K132").Value = 61002.444
$ws.Range("L132").Value = 13730.6661
$ws.Range("M132").Value = -58472.444
$ws.Range("N132").Value = -18790.6661
$ws.Range("H134").Value = 1214.9615
$ws.Range("I134").Value = 1037.125
$ws.Range("J134").Value = 1499.5
$ws.Range("K134").Value = 3111.375
$ws.Range("L134").Value = 4498.5
$ws.Range("M134").Value = -576.375
$ws.Range("N134").Value = -9568.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8399.75
$ws.Range("I2").Value = 33366.668
$ws.Range("J2").Value = 77.44444
$ws.Range("K2").Value = 200200.008
$ws.Range("L2").Value = 464.66664
$ws.Range("M2").Value = -200087.008
$ws.Range("N2").Value = -690.66664
$ws.Range("H33").Value = 176.18182
$ws.Range("I33").Value = 34.25
$ws.Range("J33").Value = 257.2857
$ws.Range("K33").Value = 205.5
$ws.Range("L33").Value = 1543.7142
$ws.Range("M33").Value = 77.5
$ws.Range("N33").Value = -2109.7142
$ws.Range("H62").Value = 6371
$ws.Range("I62").Value = 2806.5
$ws.Range("J62").Value = 8407.857
$ws.Range("K62").Value = 8419.5
$ws.Range("L62").Value = 25223.571
$ws.Range("M62").Value = -7733.5
$ws.Range("N62").Value = -26595.571
$ws.Range("H65").Value = 6371
$ws.Range("I65").Value = 2806.5
$ws.Range("J65").Value = 8407.857
$ws.Range("K65").Value = 25258.5
$ws.Range("L65").Value = 75670.713
$ws.Range("M65").Value = -21826.5
$ws.Range("N65").Value = -82534.713
$ws.Range("H68").Value = 17550.5
$ws.Range("I68").Value = 1100
$ws.Range("J68").Value = 34001
$ws.Range("K68").Value = 3300
$ws.Range("L68").Value = 102003
$ws.Range("M68").Value = -2489
$ws.Range("N68").Value = -103625
$ws.Range("H71").Value = 17550.5
$ws.Range("I71").Value = 1100
$ws.Range("J71").Value = 34001
$ws.Range("K71").Value = 9900
$ws.Range("L71").Value = 306009
$ws.Range("M71").Value = -5844
$ws.Range("N71").Value = -314121
$ws.Range("H122").Value = 720.8
$ws.Range("I122").Value = 385
$ws.Range("J122").Value = 758.1111
$ws.Range("K122").Value = 3465
$ws.Range("L122").Value = 6822.9999
$ws.Range("M122").Value = -1015
$ws.Range("N122").Value = -11722.9999
$ws.Range("H131").Value = 777.16
$ws.Range("I131").Value = 286
$ws.Range("J131").Value = 803.0105
$ws.Range("K131").Value = 858
$ws.Range("L131").Value = 2409.0315
$ws.Range("M131").Value = 4182
$ws.Range("N131").Value = -12489.0315

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15669.25
$ws.Range("I70").Value = 13070.8
$ws.Range("J70").Value = 20000
$ws.Range("K70").Value = 13070.8
$ws.Range("L70").Value = 20000
$ws.Range("M70").Value = -12800.8
$ws.Range("N70").Value = -20540
$ws.Range("H73").Value = 15669.25
$ws.Range("I73").Value = 13070.8
$ws.Range("J73").Value = 20000
$ws.Range("K73").Value = 13070.8
$ws.Range("L73").Value = 20000
$ws.Range("M73").Value = -12134.8
$ws.Range("N73").Value = -21872
$ws.Range("H97").Value = 1595.375
$ws.Range("I97").Value = 1558.6666
$ws.Range("J97").Value = 1705.5
$ws.Range("K97").Value = 1558.6666
$ws.Range("L97").Value = 1705.5
$ws.Range("M97").Value = -1062.6666
$ws.Range("N97").Value = -2697.5
$ws.Range("H132").Value = 69247.61
$ws.Range("I132").Value = 59988.723
$ws.Range("J132").Value = 102579.6
$ws.Range("K132").Value = 179966.169
$ws.Range("L132").Value = 307738.8
$ws.Range("M132").Value = -177436.169
$ws.Range("N132").Value = -312798.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 4000
$ws.Range("N4").Value = -4226
$ws.Range("H7").Value = 2990.2856
$ws.Range("I7").Value = 3190.0908
$ws.Range("J7").Value = 2770.5
$ws.Range("K7").Value = 3190.0908
$ws.Range("L7").Value = 2770.5
$ws.Range("M7").Value = -3078.0908
$ws.Range("N7").Value = -2994.5
$ws.Range("H22").Value = 1644.2106
$ws.Range("I22").Value = 1245.6666
$ws.Range("J22").Value = 2327.4285
$ws.Range("K22").Value = 1245.6666
$ws.Range("L22").Value = 2327.4285
$ws.Range("M22").Value = -950.6666
$ws.Range("N22").Value = -2917.4285
$ws.Range("H27").Value = 1644.2106
$ws.Range("I27").Value = 1245.6666
$ws.Range("J27").Value = 2327.4285
$ws.Range("K27").Value = 1245.6666
$ws.Range("L27").Value = 2327.4285
$ws.Range("M27").Value = -1138.6666
$ws.Range("N27").Value = -2541.4285
$ws.Range("H28").Value = 4000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 4000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 4000
$ws.Range("N28").Value = -4464
$ws.Range("H37").Value = 4000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 4000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 4000
$ws.Range("N37").Value = -4214
$ws.Range("H40").Value = 3530.2666
$ws.Range("I40").Value = 3625
$ws.Range("J40").Value = 3422
$ws.Range("K40").Value = 3625
$ws.Range("L40").Value = 3422
$ws.Range("M40").Value = -3489
$ws.Range("N40").Value = -3694
$ws.Range("H104").Value = 20381.428
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 20381.428
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 20381.428
$ws.Range("N104").Value = -27369.428
$ws.Range("H126").Value = 2990.2856
$ws.Range("I126").Value = 3190.0908
$ws.Range("J126").Value = 2770.5
$ws.Range("K126").Value = 9570.2724
$ws.Range("L126").Value = 8311.5
$ws.Range("M126").Value = -7100.2724
$ws.Range("N126").Value = -13251.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 127490
$ws.Range("I14").Value = 200000
$ws.Range("J14").Value = 103320
$ws.Range("K14").Value = 200000
$ws.Range("L14").Value = 103320
$ws.Range("M14").Value = -199832
$ws.Range("N14").Value = -103656
$ws.Range("H100").Value = 332.63635
$ws.Range("I100").Value = 332.63635
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 665.2727
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -124.2727
